# Generate Report for Handoff
# b.md has been handed off: update status + latest handoff file/datetime
# for the "b.md" row (row 3) on the Overview sheet and on each locale
# sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- Per-locale sheets -------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; HandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"; HandoffDate = "2016-03-08 10:21:51" },
    @{ Sheet = "de-de"; HandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"; HandoffDate = "2016-03-08 10:21:55" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Status column (B) for the b.md row (row 3)
    $ws.Range("B3").Value = $newStatus

    # Latest Handoff File column (C) for the b.md row (row 3)
    $ws.Range("C3").Value = $locale.HandoffFile

    # Keep the hyperlink on C3 pointing at the same target, but refresh
    # its displayed text to match the new handoff file name.
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq '$C$3') {
            $hl.TextToDisplay = $locale.HandoffFile
        }
    }

    # Latest Handoff Datetime column (D) for the b.md row (row 3)
    $ws.Range("D3").Value = $locale.HandoffDate
}
